# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that look like plain numbers need a leading apostrophe
# so Excel keeps them as text (matching the original inline-string cells)
# instead of silently converting them to numeric values.

$ws.Range("D2").Value = '56.550.16'
$ws.Range("E2").Value = '  -2.78%  '

$ws.Range("D3").Value = '2.984.54'
$ws.Range("E3").Value = '  -4.71%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''496.70'
$ws.Range("E5").Value = '  -5.04%  '

$ws.Range("D6").Value = '''134.96'
$ws.Range("E6").Value = '  +0.68%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -2.77%  '

$ws.Range("D9").Value = '''7.27'
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("E10").Value = '  -2.85%  '

$ws.Range("E11").Value = '  -6.20%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '3.492.25'
$ws.Range("E13").Value = '  -4.83%  '

$ws.Range("D14").Value = '''24.89'
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").Value = '56.470.89'
$ws.Range("E15").Value = '  -2.90%  '

$ws.Range("D16").Value = '2.981.10'
$ws.Range("E16").Value = '  -4.98%  '

$ws.Range("E17").Value = '  -3.43%  '

$ws.Range("D18").Value = '''5.83'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").Value = '''12.37'
$ws.Range("E19").Value = '  -4.74%  '

$ws.Range("E20").Value = '  -1.54%  '

$ws.Range("D21").Value = '''324.75'
$ws.Range("E21").Value = '  -5.01%  '

$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E23").Value = '  -7.43%  '

$ws.Range("D24").Value = '''61.28'
$ws.Range("E24").Value = '  -9.34%  '

$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("D26").Value = '''0.162'
$ws.Range("E26").Value = '  -2.28%  '

$ws.Range("E27").Value = '  -5.97%  '

$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("D29").Value = '''6.48'

$ws.Range("D30").Value = '''6.69'
$ws.Range("E30").Value = '  -1.33%  '

$ws.Range("E31").Value = '  -2.95%  '

$ws.Range("E32").Value = '  -6.88%  '

$ws.Range("D33").Value = '''20.18'
$ws.Range("E33").Value = '  -5.42%  '

$ws.Range("D34").Value = '''155.54'
$ws.Range("E34").Value = '  -0.61%  '

$ws.Range("D35").Value = '''4.46'
$ws.Range("E35").Value = '  -6.38%  '

$ws.Range("E36").Value = '  -5.78%  '

$ws.Range("D37").Value = '''5.61'
$ws.Range("E37").Value = '  -9.35%  '

$ws.Range("D38").Value = '''0.0680'
$ws.Range("E38").Value = '  -0.57%  '

$ws.Range("D39").Value = '''23.25'
$ws.Range("E39").Value = '  -3.40%  '

$ws.Range("D40").Value = '3.015.35'
$ws.Range("E40").Value = '  -4.67%  '

$ws.Range("D41").Value = '''36.56'
$ws.Range("E41").Value = '  -9.32%  '

$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").Value = '''0.637'
$ws.Range("E43").Value = '  -7.76%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.222.37'
$ws.Range("E44").Value = '  -1.48%  '

$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '''0.991'
$ws.Range("E45").Value = '  -7.94%  '

$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("E47").Value = '  -7.79%  '

$ws.Range("D48").Value = '''1.93'
$ws.Range("E48").Value = '  +5.36%  '

$ws.Range("D50").Value = '''5.79'
$ws.Range("E50").Value = '  -6.22%  '

$ws.Range("D51").Value = '''18.99'
$ws.Range("E51").Value = '  -7.49%  '
